# logboek.xlsx update — "Dit is mijn commit"
#
# Week 47: log a new activity in row 13 (11:00-11:40, "Project op github gezet")
# Week 48: fill in the first four activities of the week (rows 8-11) and
#          widen column F so the longer activity texts are readable; this
#          sheet becomes the active tab/sheet of the workbook.
# Week 1 : just a different cell selected when the file was last saved.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Week 47
# ---------------------------------------------------------------------
$week47 = $wb.Worksheets.Item("Week 47")

# New row-13 activity: 11:00 - 11:40, "Project op github gezet"
$week47.Range("C9").Copy()
$week47.Range("C13").PasteSpecial(-4122)
$week47.Range("D9").Copy()
$week47.Range("D13").PasteSpecial(-4122)
$week47.Range("C13").Value = 0.45833333333333331
$week47.Range("D13").Value = 0.4861111111111111
$week47.Range("F13").Value = "Project op github gezet"

# Rows 23-27 now use the plain date style (like B10) instead of the long
# date style that was used before.
$week47.Range("B10").Copy()
$week47.Range("B23").PasteSpecial(-4122)
$week47.Range("B24").PasteSpecial(-4122)
$week47.Range("B25").PasteSpecial(-4122)
$week47.Range("B26").PasteSpecial(-4122)
$week47.Range("B27").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Week 1 — only the remembered selection changed.
# ---------------------------------------------------------------------
$week1 = $wb.Worksheets.Item("Week 1")
$week1.Range("B35").Select()

# Selection moved off the (no longer active) tab. Done after Week 1 so
# Week 47 isn't left as the active sheet.
$week47.Range("D11:D13").Select()

# ---------------------------------------------------------------------
# Week 48 — this becomes the active sheet/tab, so it is handled last.
# ---------------------------------------------------------------------
$week48 = $wb.Worksheets.Item("Week 48")

# Rows 10 & 11 switch from the "short" time style to the centered one
# (matching rows 8 & 9) before the values are written in.
$week48.Range("C9").Copy()
$week48.Range("C10").PasteSpecial(-4122)
$week48.Range("D10").PasteSpecial(-4122)
$week48.Range("C11").PasteSpecial(-4122)
$week48.Range("D11").PasteSpecial(-4122)

$week48.Range("B8").Value = 41603
$week48.Range("C8").Value = 0.49027777777777781
$week48.Range("D8").Value = 0.50069444444444444
$week48.Range("F8").Value = "Gamenaam gewijzigd"

$week48.Range("C9").Value = 0.50069444444444444
$week48.Range("D9").Value = 0.51111111111111118
$week48.Range("F9").Value = "Hoogte en Breedte aangepast van het canvas "

$week48.Range("C10").Value = 0.51111111111111118
$week48.Range("D10").Value = 0.5180555555555556
$week48.Range("F10").Value = "icoon Toegevoegd"

$week48.Range("C11").Value = 0.51874999999999993
$week48.Range("D11").Value = 0.53125
$week48.Range("F11").Value = "Het spel laten stoppen en achtergrondkleur toegevoegd"

# Column F needs to be wide enough for the new activity descriptions.
$week48.Columns.Item(6).ColumnWidth = 52

# Week 48 becomes the active sheet/tab, scrolled down a bit.
$week48.Activate()
$week48.Application.ActiveWindow.ScrollRow = 7
$week48.Range("G28").Select()
